# Apply the edits described in the commit: "modification et suppression d'un produit analysé"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Raison sociale changed
$ws.Range("E2").Value = "jB"

# Product line (row 18) updated: quantity, designation, unit price and sub-total
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Fanta"
$ws.Range("D18").Value = 180000
$ws.Range("G18").Value = 420000
